# commit conserto dos problemas listados
#
# Applies the fixes listed in the review sheet:
#  - Rows 4 and 5: the pending "ver local" / "ver cadastro" follow-ups are
#    finished, so column D is marked "concluído" (green, centered) just like
#    the other completed rows.
#  - Row 8: no textual change (the "Bloquear campo..." text stays the same).
#  - Row 12 ("Adicionar coluna pedido cliente"): marked concluído with a
#    yellow highlight and a note in column E to "verificar folha".
#  - Rows 15 and 16: also marked concluído; row 15's component cell is
#    highlighted in red to call attention to it.
#  - Row 17: marked concluído with a yellow highlight (and green text) plus
#    the "verificar folha" note in column E.
#  - New row 18: a new pending item for the "Principal"/"Login" screen:
#    "Acertar a tela principal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colGreen  = 5287936   # RGB(0,176,80)
$colRed    = 255       # RGB(255,0,0)
$colYellow = 65535     # RGB(255,255,0)
$hCenter   = -4108     # xlCenter

# --- Row 4: "ver local" -> "concluído" ---------------------------------
$ws.Range("D4").Value = "concluído"
$ws.Range("D4").Font.Color = $colGreen
$ws.Range("D4").HorizontalAlignment = $hCenter

# --- Row 5: "ver cadastro" -> "concluído" -------------------------------
$ws.Range("D5").Value = "concluído"
$ws.Range("D5").Font.Color = $colGreen
$ws.Range("D5").HorizontalAlignment = $hCenter

# --- Row 15: highlight component cell in red ----------------------------
$ws.Range("B15").Font.Color = $colRed
$ws.Range("B15").HorizontalAlignment = $hCenter

# --- New row 18: pending item for the main/login screen -----------------
# (set first so the newly introduced shared strings are appended in the
# same order as the source workbook)
$ws.Range("C18").Value = "Acertar a tela principal"
$ws.Range("A18").Value = "Login"
$ws.Range("B18").Value = "Principal"

# --- Row 12: concluído + yellow highlight + note -------------------------
$ws.Range("D12").Value = "concluído"
$ws.Range("D12").Interior.Color = $colYellow
$ws.Range("D12").HorizontalAlignment = $hCenter

# --- Row 15 / 16: concluído -----------------------------------------------
$ws.Range("D15").Value = "concluído"
$ws.Range("D15").Font.Color = $colGreen
$ws.Range("D15").HorizontalAlignment = $hCenter

$ws.Range("D16").Value = "concluído"
$ws.Range("D16").Font.Color = $colGreen
$ws.Range("D16").HorizontalAlignment = $hCenter

# --- Row 17: concluído + yellow highlight + note --------------------------
$ws.Range("D17").Value = "concluído"
$ws.Range("D17").Interior.Color = $colYellow
$ws.Range("D17").Font.Color = $colGreen
$ws.Range("D17").HorizontalAlignment = $hCenter

# --- "verificar folha" notes in column E (rows 12 and 17) ----------------
$ws.Range("E12").Value = "verificar folha"
$ws.Range("E17").Value = "verificar folha"

# --- Selection, matching the saved view in the workbook -------------------
[void]$ws.Range("D17").Select()
